$d = $word.ActiveDocument

# --- Hunk 1: add explicit black color to the paragraph mark of the last
# paragraph ("Висновок: ...") by setting Font.Color on its full Range
# (text + paragraph mark). ---
$lastPara = $d.Paragraphs.Last
$lastPara.Range.Font.Color = 0

# --- Hunk 2: append two new paragraphs plus a blank paragraph after the
# "Висновок" paragraph, and relocate the hidden _GoBack bookmark onto the
# new final (empty) paragraph. ---

# The _GoBack bookmark currently sits, collapsed, right at the very end of
# the document (just before the final paragraph mark). Remove it so it
# does not end up attached to the wrong paragraph once we insert new
# content; we will re-create it in the correct spot afterwards.
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

$lastPara = $d.Paragraphs.Last
$insertionPoint = $d.Range($lastPara.Range.End - 1, $lastPara.Range.End - 1)

$newXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body>' +
  '<w:p>' +
    '<w:pPr>' +
      '<w:spacing w:line="240" w:lineRule="auto"/>' +
      '<w:ind w:firstLine="142"/>' +
      '<w:rPr>' +
        '<w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' +
        '<w:sz w:val="24"/>' +
        '<w:szCs w:val="24"/>' +
        '<w:lang w:eastAsia="uk-UA"/>' +
      '</w:rPr>' +
    '</w:pPr>' +
    '<w:r>' +
      '<w:rPr>' +
        '<w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' +
        '<w:sz w:val="24"/>' +
        '<w:szCs w:val="24"/>' +
        '<w:lang w:eastAsia="uk-UA"/>' +
      '</w:rPr>' +
      '<w:t>Це завдання було змінено до іншої версії, бо було виконано некоректно.</w:t>' +
    '</w:r>' +
  '</w:p>' +
  '<w:p/>' +
  '<w:p>' +
    '<w:pPr>' +
      '<w:spacing w:line="240" w:lineRule="auto"/>' +
      '<w:ind w:firstLine="142"/>' +
      '<w:rPr>' +
        '<w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' +
        '<w:sz w:val="24"/>' +
        '<w:szCs w:val="24"/>' +
        '<w:lang w:eastAsia="uk-UA"/>' +
      '</w:rPr>' +
    '</w:pPr>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
    '<w:bookmarkEnd w:id="0"/>' +
  '</w:p>' +
  '</w:body>' +
  '</w:document>' +
  '</pkg:xmlData>' +
  '</pkg:part>' +
  '</pkg:package>'

$null = $insertionPoint.InsertXML($newXml)
